$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "58dab42903e09ef89059fe914c59aa11"
$ws.Range("B11").Value = "310234a99664071f5bece0ec308a3ed8"
$ws.Range("B15").Value = "24e09ee8218b02fe5bd96aac8cfedfe4"
$ws.Range("B17").Value = "8877423d4f0b03c400ff6bed74a06713"
$ws.Range("B24").Value = "e3fdda93874f083ec3c11e1d9146626a"
$ws.Range("B29").Value = "dfa9263ed277978bbfe3b629da3fa743"
$ws.Range("B34").Value = "9365c1747950429dab505b3ab4012e69"
$ws.Range("B121").Value = "084912894f0a1eedf1f6d9b58aede52c"
$ws.Range("B133").Value = "3836913ef992e325c4850d4d8608781c"
$ws.Range("B136").Value = "4d49cffb2d9e9faa7796cbde2e091438"
$ws.Range("B159").Value = "4980632e3513d71821ba456089af32d4"
$ws.Range("B162").Value = "d26a25d453fe8716fa2abdc31a754ce7"
$ws.Range("B169").Value = "5a787603d7cbf98a13fca67b4954034d"
$ws.Range("B175").Value = "afca1fb3e1cf6e1145c2433fc8710362"
$ws.Range("B180").Value = "e6eb35d828b066bc16575e1f7c9a30ca"
$ws.Range("B183").Value = "5db62a3d8675f34a3645c5ffc0438973"
$ws.Range("B191").Value = "9925abf753c4ffda0df6660cc69afbed"
$ws.Range("B198").Value = "8584645789a5ce7a9229f2e1e71ff539"
$ws.Range("B200").Value = "79bb370ce8e3a268426f1fb339c406ef"
$ws.Range("B213").Value = "74d85772f6798c79206c67ff52853f88"
$ws.Range("B227").Value = "d4750eda165f7900d39576e86203d128"
$ws.Range("B228").Value = "4e53892efb1bce3fcff3ebdd12071fa8"
$ws.Range("B232").Value = "6bd8526ef76661fe822c47ddac716bbf"
$ws.Range("B281").Value = "9ca1e8bff17f8b6d131a92b588c4d692"
$ws.Range("B339").Value = "29687645a1a61474caf4c9376436e8a8"
$ws.Range("B460").Value = "9ac4787fb3e022cab8511c3f6f494022"
$ws.Range("B461").Value = "d35575284fa01c4af3645b4d67b29505"
$ws.Range("B478").Value = "bbccb1aa4f83544b94ad2c0721635bc1"
$ws.Range("B480").Value = "963cf9dc7ee613c0a9875594db6bc554"
$ws.Range("B500").Value = "1849f9f5dfc1d30ce1ec29f0d5b6bf80"
$ws.Range("B501").Value = "7a175621f2199ecd3b1499a9fa806df8"
$ws.Range("B502").Value = "ced02ad428a2dfa62a2da9d740ebcb65"
$ws.Range("B506").Value = "6b8810479b19724e23c5fbbcce990ef0"
$ws.Range("B514").Value = "932097ea6bb4489c5b9280833f3a36a2"
$ws.Range("B517").Value = "4be45ec5d7203c0f203c534372ffd6b2"
$ws.Range("B524").Value = "840c2f9ca3e3894a78e755b0a48ca247"
$ws.Range("B547").Value = "d035b556281d3c22005bebe8f29c8f25"
$ws.Range("B563").Value = "26447e076f2e40f7dc608bbd3b9a32e8"
$ws.Range("B572").Value = "5e85027f5af3e96405dd27baa6ff8ea7"
$ws.Range("B616").Value = "a386210b38695e0815d8d0f3785aa61a"
$ws.Range("B627").Value = "381c458e72d6669e6ebd55a2872754ae"
$ws.Range("B629").Value = "621f1a624f0c177940632e0b66545362"
$ws.Range("B649").Value = "c6bd531f4396c0d299c990fa5332263b"
$ws.Range("B655").Value = "dce1d1a5f52828bd78aacac09d4acd5a"
$ws.Range("B665").Value = "ef35b58cc3575c62094ae180e6f2d274"
$ws.Range("B666").Value = "5c18423e1edd26a7ee761fe54ba4c44a"
$ws.Range("B680").Value = "2d1aac6ae477f530a7a6ca9df04660bf"
$ws.Range("B685").Value = "25596f783d9ee2d76ca7d4f911f4cb2e"
$ws.Range("B700").Value = "2dfd3c6498de42ccecbd4ac4a5dddebb"
$ws.Range("B703").Value = "6b4b64e80de6aabc05133a0a342aba0e"
$ws.Range("B704").Value = "4a8f6aa0dfde24d93bb20b4f9f137089"
$ws.Range("B715").Value = "0fda5def2e1a945a787dffc643ffd41c"
$ws.Range("B729").Value = "e9149633f613f3ec2035f7984aaf342d"
$ws.Range("B733").Value = "c7fc4a56ff04fe1256e076b965dcb380"
$ws.Range("B742").Value = "2ba01287e3d1d91cdea937690dfceb5f"
$ws.Range("B819").Value = "a238374370b3c4710094c915ac49dccb"
$ws.Range("B830").Value = "a0dda3d45be58bcf6571bb57f0e7f349"
$ws.Range("B835").Value = "f33bd6de12419e7a4ba94cc4aca981b2"
$ws.Range("B854").Value = "7a0ca750b77933985eaab180181cb133"
